$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds free-text price strings (dotted thousands separators, many
# trailing zeros) that must stay text, not be auto-coerced to numbers/dates by
# COM assignment. Force text format on the whole D range first, write the new
# values, then restore the default style so no stray number format lingers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "29.067.71"
$ws.Range("E2").Value = "  -0.20%  "

# Row 3
$ws.Range("D3").Value = "1.820.09"
$ws.Range("E3").Value = "  -0.78%  "

# Row 4
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").Value = "241.07"
$ws.Range("E5").Value = "  -1.06%  "

# Row 6
$ws.Range("D6").Value = "0.6143"
$ws.Range("E6").Value = "  -2.30%  "

# Row 7
$ws.Range("D7").Value = "0.9986"
$ws.Range("E7").Value = "  -0.44%  "

# Row 8
$ws.Range("D8").Value = "0.07318"
$ws.Range("E8").Value = "  -2.32%  "

# Row 9
$ws.Range("D9").Value = "0.2877"
$ws.Range("E9").Value = "  -1.67%  "

# Row 10
$ws.Range("D10").Value = "22.88"
$ws.Range("E10").Value = "  -1.54%  "

# Row 11
$ws.Range("D11").Value = "0.07650"
$ws.Range("E11").Value = "  -0.53%  "

# Row 12
$ws.Range("D12").Value = "1.815.96"
$ws.Range("E12").Value = "  -0.73%  "

# Row 13
$ws.Range("D13").Value = "4.940"
$ws.Range("E13").Value = "  -1.41%  "

# Row 14
$ws.Range("D14").Value = "0.6578"
$ws.Range("E14").Value = "  -1.52%  "

# Row 15
$ws.Range("D15").Value = "81.62"
$ws.Range("E15").Value = "  -1.42%  "

# Row 16
$ws.Range("D16").Value = "0.000008944"
$ws.Range("E16").Value = "  -4.39%  "

# Row 17
$ws.Range("D17").Value = "5.813"
$ws.Range("E17").Value = "  -2.88%  "

# Row 18
$ws.Range("D18").Value = "29.042.31"
$ws.Range("E18").Value = "  -0.26%  "

# Row 19
$ws.Range("D19").Value = "2.067.08"
$ws.Range("E19").Value = "  -0.55%  "

# Row 20
$ws.Range("D20").Value = "237.30"
$ws.Range("E20").Value = "  +6.27%  "

# Row 21
$ws.Range("D21").Value = "12.40"
$ws.Range("E21").Value = "  -1.43%  "

# Row 22
$ws.Range("D22").Value = "0.9984"
$ws.Range("E22").Value = "  -0.57%  "

# Row 23
$ws.Range("D23").Value = "7.098"
$ws.Range("E23").Value = "  -0.07%  "

# Row 24
$ws.Range("D24").Value = "0.9994"
$ws.Range("E24").Value = "  -0.26%  "

# Row 25
$ws.Range("D25").Value = "157.00"
$ws.Range("E25").Value = "  -1.93%  "

# Row 26
$ws.Range("D26").Value = "0.1407"
$ws.Range("E26").Value = "  +1.16%  "

# Row 27
$ws.Range("D27").Value = "8.410"
$ws.Range("E27").Value = "  -0.96%  "

# Row 28
$ws.Range("D28").Value = "17.58"
$ws.Range("E28").Value = "  -1.74%  "

# Row 29
$ws.Range("D29").Value = "1.483"
$ws.Range("E29").Value = "  -1.16%  "

# Row 30
$ws.Range("D30").Value = "0.05539"
$ws.Range("E30").Value = "  -2.93%  "

# Row 31
$ws.Range("D31").Value = "4.080"
$ws.Range("E31").Value = "  +0.02%  "

# Row 32
$ws.Range("D32").Value = "4.079"
$ws.Range("E32").Value = "  -1.74%  "

# Row 33
$ws.Range("D33").Value = "1.200"
$ws.Range("E33").Value = "  -0.66%  "

# Row 34
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "1.817"
$ws.Range("E34").Value = "  -0.82%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.7328"
$ws.Range("E35").Value = "  -1.18%  "

# Row 36
$ws.Range("D36").Value = "1.127"
$ws.Range("E36").Value = "  -1.14%  "

# Row 37
$ws.Range("E37").Value = "  -2.63%  "

# Row 38
$ws.Range("E38").Value = "  +2.35%  "

# Row 39
$ws.Range("D39").Value = "1.206.60"
$ws.Range("E39").Value = "  -0.65%  "

# Row 40
$ws.Range("D40").Value = "0.01753"
$ws.Range("E40").Value = "  -1.44%  "

# Row 41
$ws.Range("E41").Value = "  -2.63%  "

# Row 42
$ws.Range("D42").Value = "0.8940"
$ws.Range("E42").Value = "  +0.56%  "

# Row 43
$ws.Range("D43").Value = "0.9981"
$ws.Range("E43").Value = "  -0.48%  "

# Row 44
$ws.Range("D44").Value = "100.98"
$ws.Range("E44").Value = "  -0.97%  "

# Row 45
$ws.Range("D45").Value = "1.964.04"
$ws.Range("E45").Value = "  -0.88%  "

# Row 46
$ws.Range("D46").Value = "64.39"
$ws.Range("E46").Value = "  -1.82%  "

# Row 47
$ws.Range("E47").Value = "  -0.23%  "

# Row 48
$ws.Range("D48").Value = "0.00000000119"
$ws.Range("E48").Value = "  -4.88%  "

# Row 49
$ws.Range("D49").Value = "0.3988"
$ws.Range("E49").Value = "  -1.85%  "

# Row 50
$ws.Range("D50").Value = "8.992"
$ws.Range("E50").Value = "  -0.07%  "

# Row 51
$ws.Range("D51").Value = "0.05748"
$ws.Range("E51").Value = "  -1.21%  "

# Restore the default (General/Normal) style on column D now that the values
# are committed as text, so the saved styles table matches the original shape.
$dRange.Style = "Normal"
